$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay text even when it looks numeric
# (mirrors Excel's own behaviour for cells pre-formatted as Text), then
# drop the temporary Text number-format again so no stray style sticks.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

# Row 2: fewer coexisting ions, new frame count
$ws.Range("A2").Value = "130, 786"
Set-TextValue "D2" "959"

# Row 3: single coexisting ion, permeated residue + frame count change
Set-TextValue "A3" "130"
Set-TextValue "C3" "130"
Set-TextValue "D3" "969"

# Row 4: new residue combination entirely
$ws.Range("A4").Value = "98, 1082, SF"
$ws.Range("C4").Value = "SF"
Set-TextValue "D4" "806"

# Row 5: new residue combination entirely
$ws.Range("A5").Value = "754, 786, 1082"
Set-TextValue "C5" "1082"
Set-TextValue "D5" "946"

# Rows 6-16 no longer exist in the updated results - remove them entirely
$ws.Range("A6:A16").EntireRow.Delete()
